$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7 and 8 hold two different players' full stat lines (columns A:DL).
# This edit swaps the two rows (Jorgen Strand Larsen moves from row 8 to row 7,
# Christantus Uche moves from row 7 to row 8), and also updates the
# totwAppearances figure (column I) for Jorgen Strand Larsen's row from 0 to 1
# as part of the refreshed Betting Markets Analytics data.

$row7vals = $ws.Range("A7:DL7").Value2
$row8vals = $ws.Range("A8:DL8").Value2

$ws.Range("A7:DL7").Value2 = $row8vals
$ws.Range("A8:DL8").Value2 = $row7vals

# totwAppearances (column I) for the player now in row 7 (Jorgen Strand Larsen)
$ws.Range("I7").Value2 = 1
